# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback
# timestamps for the 82752071-1c07-4ffe-91f3-7a0b4128ebe5.md entry
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H4").Value = "2016-08-13 10:56:12"
$wsZhCn.Range("K4").Value = "2016-08-13 10:56:42"

# de-de: Correspond Handback DateTime
$wsDeDe.Range("K4").Value = "2016-08-13 10:56:51"

# de-de: Correspond Handoff Datetime and Overview: Latest HO Xliff Generate Date
# both held the same timestamp before the edit and both move to the same
# new timestamp now.
$wsDeDe.Range("H4").Value = "2016-08-13 10:56:19"
$wsOverview.Range("G4").Value = "2016-08-13 10:56:19"
